# Auto-generated edit script applying numeric updates to Kraken_Profits workbook
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 198.28572
$ws.Range("I2").Value = 147.6
$ws.Range("K2").Value = 147.6
$ws.Range("M2").Value = -34.59999999999999
$ws.Range("H15").Value = 86.94444
$ws.Range("I15").Value = 86.94444
$ws.Range("K15").Value = 260.83332
$ws.Range("M15").Value = -91.83332000000001
$ws.Range("H32").Value = 8999.429
$ws.Range("J32").Value = 8999.429
$ws.Range("L32").Value = 8999.429
$ws.Range("N32").Value = -9651.429
$ws.Range("H53").Value = 314.15384
$ws.Range("I53").Value = 371.5
$ws.Range("K53").Value = 371.5
$ws.Range("M53").Value = 265.5
$ws.Range("H98").Value = 2125.4443
$ws.Range("I98").Value = 1834.875
$ws.Range("K98").Value = 1834.875
$ws.Range("M98").Value = -336.875
$ws.Range("H122").Value = 2125.4443
$ws.Range("I122").Value = 1834.875
$ws.Range("K122").Value = 5504.625
$ws.Range("M122").Value = -3054.625
$ws.Range("H137").Value = 2500
$ws.Range("I137").Value = 2500
$ws.Range("K137").Value = 7500
$ws.Range("M137").Value = -4950

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 537.6667
$ws.Range("J2").Value = 556.5
$ws.Range("L2").Value = 556.5
$ws.Range("N2").Value = -782.5
$ws.Range("H110").Value = 981.25
$ws.Range("I110").Value = 981.25
$ws.Range("K110").Value = 981.25
$ws.Range("M110").Value = 1063.75
$ws.Range("H116").Value = 537.6667
$ws.Range("J116").Value = 556.5
$ws.Range("L116").Value = 556.5
$ws.Range("N116").Value = -5144.5
$ws.Range("H122").Value = 2000
$ws.Range("I122").Value = 1500
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 4500
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -2050
$ws.Range("N122").Value = -12400

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 537.6667
$ws.Range("J3").Value = 556.5
$ws.Range("L3").Value = 556.5
$ws.Range("N3").Value = -784.5
$ws.Range("H107").Value = 800
$ws.Range("I107").Value = 600
$ws.Range("K107").Value = 600
$ws.Range("M107").Value = 1320
$ws.Range("H134").Value = 3587.4
$ws.Range("I134").Value = 3480.75
$ws.Range("K134").Value = 10442.25
$ws.Range("M134").Value = -7907.25

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("H31").Value = 2874.2727
$ws.Range("I31").Value = 1435.3077
$ws.Range("K31").Value = 1435.3077
$ws.Range("M31").Value = -1140.3077
$ws.Range("H34").Value = 2874.2727
$ws.Range("I34").Value = 1435.3077
$ws.Range("K34").Value = 1435.3077
$ws.Range("M34").Value = -1233.3077
$ws.Range("H82").Value = 59000
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 59000
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 100.4
$ws.Range("I26").Value = 25.5
$ws.Range("J26").Value = 150.33333
$ws.Range("K26").Value = 76.5
$ws.Range("L26").Value = 450.99999
$ws.Range("M26").Value = 211.5
$ws.Range("N26").Value = -1026.99999
$ws.Range("H107").Value = 383
$ws.Range("I107").Value = 74.5
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 223.5
$ws.Range("L107").Value = 3000
$ws.Range("M107").Value = 1696.5
$ws.Range("N107").Value = -6840

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 213.66667
$ws.Range("I107").Value = 213.66667
$ws.Range("K107").Value = 213.66667
$ws.Range("M107").Value = 1706.33333
$ws.Range("H113").Value = 2194
$ws.Range("J113").Value = 2888
$ws.Range("L113").Value = 2888
$ws.Range("N113").Value = -7228

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7123.75
$ws.Range("I7").Value = 6831.6665
$ws.Range("K7").Value = 6831.6665
$ws.Range("M7").Value = -6719.6665
$ws.Range("H22").Value = 1237.1111
$ws.Range("I22").Value = 1646.8
$ws.Range("J22").Value = 725
$ws.Range("K22").Value = 1646.8
$ws.Range("L22").Value = 725
$ws.Range("M22").Value = -1351.8
$ws.Range("N22").Value = -1315
$ws.Range("H27").Value = 1237.1111
$ws.Range("I27").Value = 1646.8
$ws.Range("J27").Value = 725
$ws.Range("K27").Value = 1646.8
$ws.Range("L27").Value = 725
$ws.Range("M27").Value = -1539.8
$ws.Range("N27").Value = -939
$ws.Range("H61").Value = 4995.75
$ws.Range("I61").Value = 5666.3335
$ws.Range("K61").Value = 5666.3335
$ws.Range("M61").Value = -5464.3335
$ws.Range("H93").Value = 1900.5714
$ws.Range("I93").Value = 1867.3334
$ws.Range("K93").Value = 1867.3334
$ws.Range("M93").Value = -619.3334
$ws.Range("H113").Value = 4995.75
$ws.Range("I113").Value = 5666.3335
$ws.Range("K113").Value = 5666.3335
$ws.Range("M113").Value = -3496.3335
$ws.Range("H122").Value = 4599.75
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()
$ws.Range("H126").Value = 7123.75
$ws.Range("I126").Value = 6831.6665
$ws.Range("K126").Value = 20494.9995
$ws.Range("M126").Value = -18024.9995

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 50001
$ws.Range("I81").Value = 50001
$ws.Range("K81").Value = 100002
$ws.Range("M81").Value = -98941
$ws.Range("H84").Value = 50001
$ws.Range("I84").Value = 50001
$ws.Range("K84").Value = 500010
$ws.Range("M84").Value = -494706
$ws.Range("H113").Value = 1005.1
$ws.Range("J113").Value = 774.5
$ws.Range("L113").Value = 2323.5
$ws.Range("N113").Value = -6663.5
$ws.Range("H126").Value = 4272.778
$ws.Range("I126").Value = 3076
$ws.Range("K126").Value = 9228
$ws.Range("M126").Value = -6758
$ws.Range("H132").Value = 4500
$ws.Range("I132").Value = 4000
$ws.Range("K132").Value = 12000
$ws.Range("M132").Value = -9470
$ws.Range("H136").Value = 6499.8
$ws.Range("I136").Value = 6124.75
$ws.Range("K136").Value = 18374.25
$ws.Range("M136").Value = -15824.25
